$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.65
$ws.Range("G2").Value = 1.67
$ws.Range("H2").Value = 5.5
$ws.Range("I2").Value = 5.7
$ws.Range("J2").Value = 4.6
$ws.Range("K2").Value = 4.8
$ws.Range("N2").Value = 4.8
$ws.Range("Q2").Value = 1.69
$ws.Range("S2").Value = 2.72
$ws.Range("U2").Value = 2.16
$ws.Range("V2").Value = 1.21
$ws.Range("W2").Value = 2.5
$ws.Range("X2").Value = 26
$ws.Range("Y2").Value = 25
$ws.Range("AD2").Value = 25
$ws.Range("AE2").Value = 85
$ws.Range("AH2").Value = 22
$ws.Range("AL2").Value = 34
$ws.Range("AN2").Value = 8

$ws.Range("G3").Value = 2.6
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.12
$ws.Range("N3").Value = 2.34
$ws.Range("O3").Value = 1.6
$ws.Range("Q3").Value = 2.6
$ws.Range("R3").Value = 1.15
$ws.Range("S3").Value = 5.4
$ws.Range("T3").Value = 2.04
$ws.Range("U3").Value = 1.56
$ws.Range("V3").Value = 1.24
$ws.Range("W3").Value = 1.63
$ws.Range("X3").Value = 980
$ws.Range("Y3").Value = 13
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 7.8
$ws.Range("AC3").Value = 8.6
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

$ws.Range("J5").Value = 2.56
$ws.Range("P5").Value = 1.41
$ws.Range("Q5").Value = 2.7

$ws.Range("H6").Value = 5.1
$ws.Range("I6").Value = 16.5
$ws.Range("J6").Value = 3.95
$ws.Range("K6").Value = 8.6
$ws.Range("P6").Value = 1.99

$ws.Range("AJ9").Value = 19

$ws.Range("F10").Value = 1.12
$ws.Range("G10").Value = 1.14
$ws.Range("H10").Value = 25
$ws.Range("I10").Value = 95
$ws.Range("J10").Value = 10.5
$ws.Range("K10").Value = 12
